$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 2.55
$ws.Range("I3").Value = 3.25
$ws.Range("J3").Value = 3.4
$ws.Range("O3").Value = 1.62
$ws.Range("P3").Value = 2.2
$ws.Range("U3").Value = 5.2
$ws.Range("V3").Value = 1.17
$ws.Range("W3").Value = 6.5
$ws.Range("X3").Value = 1.11
$ws.Range("AD3").Value = 10
$ws.Range("AG3").Value = 26
$ws.Range("AN3").Value = 7
$ws.Range("G6").Value = 2.3
$ws.Range("I6").Value = 3
$ws.Range("K6").Value = 2.2
$ws.Range("L6").Value = 3.5
$ws.Range("M6").Value = 1.05
$ws.Range("N6").Value = 11
$ws.Range("AA6").Value = 1.67
$ws.Range("AB6").Value = 2.1
$ws.Range("AC6").Value = 9
$ws.Range("AE6").Value = 9
$ws.Range("AF6").Value = 21
$ws.Range("AJ6").Value = 7
$ws.Range("AK6").Value = 13
$ws.Range("AO6").Value = 15
$ws.Range("AR6").Value = 23
$ws.Range("H9").Value = 3.95
$ws.Range("I9").Value = 8.75
$ws.Range("L9").Value = 7.6
$ws.Range("X9").Value = 1.28
$ws.Range("AA9").Value = 2.25
$ws.Range("AB9").Value = 1.57
$ws.Range("AD9").Value = 5.7
$ws.Range("AF9").Value = 8.75
$ws.Range("AG9").Value = 13
$ws.Range("AN9").Value = 17.5
$ws.Range("G10").Value = 1.75
$ws.Range("H10").Value = 3.7
$ws.Range("I10").Value = 3.8
$ws.Range("J10").Value = 2.38
$ws.Range("AJ10").Value = 7.5
$ws.Range("AN10").Value = 15
$ws.Range("AP10").Value = 13
$ws.Range("G11").Value = 1.9
$ws.Range("H11").Value = 3.4
$ws.Range("I11").Value = 4.2
$ws.Range("M11").Value = 1.05
$ws.Range("N11").Value = 11
$ws.Range("O11").Value = 1.25
$ws.Range("P11").Value = 3.75
$ws.Range("S11").Value = 1.85
$ws.Range("W11").Value = 3
$ws.Range("X11").Value = 1.36
$ws.Range("AA11").Value = 1.7
$ws.Range("AB11").Value = 2.05
$ws.Range("AC11").Value = 8
$ws.Range("AD11").Value = 9.5
$ws.Range("AI11").Value = 11
$ws.Range("AR11").Value = 34
$ws.Range("G12").Value = 1.22
$ws.Range("H12").Value = 5.75
$ws.Range("L12").Value = 8.5
$ws.Range("U12").Value = 1.9
$ws.Range("V12").Value = 1.9
$ws.Range("AE12").Value = 10
$ws.Range("AG12").Value = 12
$ws.Range("AJ12").Value = 12
$ws.Range("AK12").Value = 26
$ws.Range("AP12").Value = 26
$ws.Range("G13").Value = 1.67
$ws.Range("I13").Value = 4.1
$ws.Range("L13").Value = 4.5
$ws.Range("S13").Value = 1.65
$ws.Range("T13").Value = 2.2
$ws.Range("Y13").Value = 1.33
$ws.Range("Z13").Value = 3.25
$ws.Range("AA13").Value = 1.73
$ws.Range("AB13").Value = 2
$ws.Range("AC13").Value = 8.5
$ws.Range("AH13").Value = 23
$ws.Range("AP13").Value = 15
$ws.Range("AR13").Value = 34
$ws.Range("M14").Value = 1.02
$ws.Range("N14").Value = 11
$ws.Range("O14").Value = 1.29
$ws.Range("P14").Value = 3.5
$ws.Range("G15").Value = 1.67
$ws.Range("S15").Value = 1.48
$ws.Range("I16").Value = 2
$ws.Range("G17").Value = 5.7
$ws.Range("H17").Value = 4.25
$ws.Range("I17").Value = 1.5
$ws.Range("J17").Value = 5.5
$ws.Range("K17").Value = 2.37
$ws.Range("L17").Value = 1.98
$ws.Range("N17").Value = 8.5
$ws.Range("O17").Value = 1.21
$ws.Range("P17").Value = 3.9
$ws.Range("S17").Value = 1.65
$ws.Range("T17").Value = 2.12
$ws.Range("W17").Value = 2.55
$ws.Range("X17").Value = 1.45
$ws.Range("Y17").Value = 1.32
$ws.Range("Z17").Value = 3.1
$ws.Range("AA17").Value = 1.78
$ws.Range("AB17").Value = 1.93
$ws.Range("AC17").Value = 17.5
$ws.Range("AD17").Value = 37
$ws.Range("AE17").Value = 18
$ws.Range("AF17").Value = 110
$ws.Range("AG17").Value = 55
$ws.Range("AI17").Value = 8.5
$ws.Range("AJ17").Value = 8.25
$ws.Range("AK17").Value = 16.5
$ws.Range("AN17").Value = 7.6
$ws.Range("AO17").Value = 7.4
$ws.Range("AP17").Value = 8
$ws.Range("AQ17").Value = 10.5
$ws.Range("AR17").Value = 11.5
$ws.Range("G18").Value = 2.07
$ws.Range("H18").Value = 3.3
$ws.Range("I18").Value = 3.3
$ws.Range("J18").Value = 2.65
$ws.Range("K18").Value = 2.07
$ws.Range("L18").Value = 3.8
$ws.Range("O18").Value = 1.36
$ws.Range("P18").Value = 2.67
$ws.Range("S18").Value = 2.05
$ws.Range("T18").Value = 1.6
$ws.Range("W18").Value = 3.4
$ws.Range("X18").Value = 1.22
$ws.Range("Y18").Value = 1.4
$ws.Range("Z18").Value = 2.5
$ws.Range("AA18").Value = 1.87
$ws.Range("AB18").Value = 1.75
$ws.Range("AC18").Value = 6.6
$ws.Range("AD18").Value = 9.25
$ws.Range("AE18").Value = 9
$ws.Range("AF18").Value = 18.5
$ws.Range("AG18").Value = 18
$ws.Range("AH18").Value = 35
$ws.Range("AI18").Value = 8.25
$ws.Range("AJ18").Value = 6.4
$ws.Range("AK18").Value = 16.5
$ws.Range("AN18").Value = 8.75
$ws.Range("AO18").Value = 16
$ws.Range("AP18").Value = 12
$ws.Range("AQ18").Value = 45
$ws.Range("AR18").Value = 32
$ws.Range("AS18").Value = 45
